$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.718.40"
$ws.Range("E2").Value = "'  -1.89%  "

$ws.Range("D3").Value = "'1.894.18"
$ws.Range("E3").Value = "'  -1.30%  "

$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "'  -0.59%  "

$ws.Range("D5").Value = "'311.55"
$ws.Range("E5").Value = "'  -1.47%  "

$ws.Range("E6").Value = "'  -0.51%  "

$ws.Range("D7").Value = "'0.4939"
$ws.Range("E7").Value = "'  +1.66%  "

$ws.Range("D8").Value = "'0.3790"
$ws.Range("E8").Value = "'  -1.59%  "

$ws.Range("D9").Value = "'0.07325"
$ws.Range("E9").Value = "'  -1.09%  "

$ws.Range("D10").Value = "'0.9096"
$ws.Range("E10").Value = "'  -4.58%  "

$ws.Range("E11").Value = "'  -2.19%  "

$ws.Range("D12").Value = "'0.07630"
$ws.Range("E12").Value = "'  -2.35%  "

$ws.Range("D13").Value = "'1.901.76"
$ws.Range("E13").Value = "'  -1.18%  "

$ws.Range("D14").Value = "'5.466"

$ws.Range("D15").Value = "'6.656"
$ws.Range("E15").Value = "'  -0.45%  "

$ws.Range("D16").Value = "'91.19"
$ws.Range("E16").Value = "'  -1.25%  "

$ws.Range("D17").Value = "'1.001"
$ws.Range("E17").Value = "'  -0.58%  "

$ws.Range("D18").Value = "'0.000008741"
$ws.Range("E18").Value = "'  -1.95%  "

$ws.Range("D19").Value = "'1.000"
$ws.Range("E19").Value = "'  -0.56%  "

$ws.Range("D20").Value = "'27.733.41"
$ws.Range("E20").Value = "'  -1.86%  "

$ws.Range("D21").Value = "'14.48"
$ws.Range("E21").Value = "'  -3.98%  "

$ws.Range("D22").Value = "'5.122"
$ws.Range("E22").Value = "'  -1.11%  "

$ws.Range("D23").Value = "'2.120.41"
$ws.Range("E23").Value = "'  -1.98%  "

$ws.Range("D24").Value = "'10.75"
$ws.Range("E24").Value = "'  -1.61%  "

$ws.Range("D25").Value = "'154.06"
$ws.Range("E25").Value = "'  -1.55%  "

$ws.Range("D26").Value = "'1.847"
$ws.Range("E26").Value = "'  -4.50%  "

$ws.Range("D27").Value = "'18.41"
$ws.Range("E27").Value = "'  -1.08%  "

$ws.Range("D28").Value = "'2.168"
$ws.Range("E28").Value = "'  +1.79%  "

$ws.Range("D29").Value = "'115.33"
$ws.Range("E29").Value = "'  -1.59%  "

$ws.Range("D30").Value = "'4.877"
$ws.Range("E30").Value = "'  -3.49%  "

$ws.Range("D31").Value = "'0.08938"
$ws.Range("E31").Value = "'  +0.15%  "

$ws.Range("D32").Value = "'3.211"
$ws.Range("E32").Value = "'  -4.54%  "

$ws.Range("D33").Value = "'1.227"
$ws.Range("E33").Value = "'  -2.18%  "

$ws.Range("D34").Value = "'0.7664"
$ws.Range("E34").Value = "'  -2.46%  "

$ws.Range("E35").Value = "'  -1.30%  "

$ws.Range("D36").Value = "'2.563"
$ws.Range("E36").Value = "'  -7.94%  "

$ws.Range("D37").Value = "'0.02043"
$ws.Range("E37").Value = "'  -0.83%  "

$ws.Range("E38").Value = "'  -3.33%  "

$ws.Range("B39").Value = "'TheSandbox"
$ws.Range("C39").Value = "'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D39").Value = "'0.5504"
$ws.Range("E39").Value = "'  -1.50%  "

$ws.Range("B40").Value = "'Hedera"
$ws.Range("C40").Value = "'https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D40").Value = "'0.05289"
$ws.Range("E40").Value = "'  -2.01%  "

$ws.Range("E41").Value = "'  -1.48%  "

$ws.Range("D42").Value = "'6.888"
$ws.Range("E42").Value = "'  -3.75%  "

$ws.Range("D43").Value = "'8.562"
$ws.Range("E43").Value = "'  -0.87%  "

$ws.Range("D44").Value = "'112.44"
$ws.Range("E44").Value = "'  +3.97%  "

$ws.Range("D45").Value = "'0.1518"
$ws.Range("E45").Value = "'  -1.53%  "

$ws.Range("D46").Value = "'10.60"
$ws.Range("E46").Value = "'  -2.65%  "

$ws.Range("D47").Value = "'0.4801"
$ws.Range("E47").Value = "'  -3.15%  "

$ws.Range("D48").Value = "'1.001"
$ws.Range("E48").Value = "'  -0.52%  "

$ws.Range("E49").Value = "'  -2.83%  "

$ws.Range("D50").Value = "'67.43"
$ws.Range("E50").Value = "'  -3.02%  "

$ws.Range("D51").Value = "'0.06050"
$ws.Range("E51").Value = "'  -1.69%  "

